$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F1").Value = "Total 12month unduplicated headcount"
$ws.Range("G1").Value = "ug unduplicated headcount"
$ws.Range("H1").Value = "g unduplicated headcount"
$ws.Range("I1").Value = "Total 12month fte enrollment"
$ws.Range("J1").Value = "ug FTE"
$ws.Range("K1").Value = "g FTE"

$ws.Range("K1").Select()
